$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp message in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 06:46"

# Apply per-cell country reorder + updated case numbers
$ws.Range("A43").Value = "India"
$ws.Range("B43").Value = 511
$ws.Range("C43").Value = 12
$ws.Range("D43").Value = 37
$ws.Range("E43").Value = 464
$ws.Range("F43").Value = 0
$ws.Range("H43").Value = 10
$ws.Range("A44").Value = "Singapur"
$ws.Range("B44").Value = 509
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 152
$ws.Range("E44").Value = 355
$ws.Range("F44").Value = 14
$ws.Range("H44").Value = 2
$ws.Range("B100").Value = 67
$ws.Range("C100").Value = 5
$ws.Range("E100").Value = 67
$ws.Range("B109").Value = 49
$ws.Range("C109").Value = 3
$ws.Range("E109").Value = 49
$ws.Range("A114").Value = "Consejo Danes para los Refugiados"
$ws.Range("F114").Value = 0
$ws.Range("A115").Value = "Mauricio"
$ws.Range("F115").Value = 1
$ws.Range("A116").Value = "Guam"
$ws.Range("C116").Value = 4
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 32
$ws.Range("H116").Value = 1
$ws.Range("A117").Value = "Banglades"
$ws.Range("B117").Value = 33
$ws.Range("D117").Value = 5
$ws.Range("E117").Value = 25
$ws.Range("H117").Value = 3
$ws.Range("A118").Value = "Puerto Rico"
$ws.Range("B118").Value = 31
$ws.Range("D118").Value = 1
$ws.Range("E118").Value = 28
$ws.Range("H118").Value = 2
$ws.Range("A119").Value = "Honduras"
$ws.Range("B119").Value = 30
$ws.Range("E119").Value = 30
$ws.Range("H119").Value = 0
$ws.Range("B120").Value = 28
$ws.Range("C120").Value = 1
$ws.Range("E120").Value = 28
$ws.Range("A122").Value = "Paraguay"
$ws.Range("C122").Value = 5
$ws.Range("F122").Value = 1
$ws.Range("G122").Value = 1
$ws.Range("A123").Value = "Ghana"
$ws.Range("C123").Value = 0
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("A129").Value = "Guyana"
$ws.Range("A130").Value = "Guatemala"
$ws.Range("A135").Value = "Madagascar"
$ws.Range("C135").Value = 5
$ws.Range("A136").Value = "Barbados"
$ws.Range("B136").Value = 17
$ws.Range("E136").Value = 17
$ws.Range("A138").Value = "Kenia"
$ws.Range("B138").Value = 16
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 16
$ws.Range("A139").Value = "Gibraltar"
$ws.Range("B139").Value = 15
$ws.Range("D139").Value = 5
$ws.Range("E139").Value = 10
$ws.Range("A140").Value = "Isla de Man"
$ws.Range("D140").Value = 0
$ws.Range("E140").Value = 13
$ws.Range("A141").Value = "Maldivas"
$ws.Range("B141").Value = 13
$ws.Range("D141").Value = 5
$ws.Range("E141").Value = 8
$ws.Range("A142").Value = "Tanzania"
$ws.Range("A146").Value = "Uganda"
$ws.Range("A147").Value = "Guinea Ecuatorial"
$ws.Range("A148").Value = "Nueva Caledonia"
$ws.Range("A149").Value = "San Martin (Parte Francesa)"
$ws.Range("A151").Value = "Benin"
$ws.Range("C151").Value = 0
$ws.Range("A153").Value = "Surinam"
$ws.Range("C153").Value = 1
$ws.Range("A158").Value = "Fiyi"
$ws.Range("C158").Value = 1
$ws.Range("A159").Value = "Guinea"
$ws.Range("A161").Value = "Congo"
$ws.Range("C161").Value = 0
$ws.Range("A162").Value = "Suazilandia"
$ws.Range("A177").Value = "San Martin (Parte Holandesa)"
$ws.Range("A179").Value = "Dominica"
$ws.Range("A180").Value = "Mauritania"
$ws.Range("A181").Value = "Butan"
$ws.Range("A182").Value = "Republica del Chad"
$ws.Range("A183").Value = "Birmania"
$ws.Range("A184").Value = "Sudan"
$ws.Range("A185").Value = "Gambia"
$ws.Range("D185").Value = 0
$ws.Range("H185").Value = 1
$ws.Range("A186").Value = "Nepal"
$ws.Range("D186").Value = 1
$ws.Range("H186").Value = 0
$ws.Range("A188").Value = "Papua Nueva Guinea"
$ws.Range("A189").Value = "Montserrat"
$ws.Range("A190").Value = "Eritrea"
$ws.Range("A191").Value = "Santa Sede"
$ws.Range("A192").Value = "Islas Turcas y Caicos"
$ws.Range("A193").Value = "Mozambique"
$ws.Range("A194").Value = "Siria"
$ws.Range("A195").Value = "Granada"
$ws.Range("A196").Value = "Somalia"
$ws.Range("A197").Value = "Belice"
$ws.Range("A198").Value = "Timor Oriental"
